$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.67
$ws.Range("H3").Value = 3.5
$ws.Range("I3").Value = 6
$ws.Range("AA3").Value = 2.5
$ws.Range("AB3").Value = 1.5
$ws.Range("AN3").Value = 10
$ws.Range("AO3").Value = 26

# Row 4
$ws.Range("G4").Value = 1.6
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 6.25
$ws.Range("AC4").Value = 5.5

# Row 5
$ws.Range("S5").Value = 2.3
$ws.Range("T5").Value = 1.6

# Row 8 (previously empty inline strings, now numeric)
$ws.Range("Q8").Value = 1.64
$ws.Range("R8").Value = 2.21
$ws.Range("U8").Value = 3.25
$ws.Range("V8").Value = 1.33

# Row 13
$ws.Range("G13").Value = 1.75
$ws.Range("I13").Value = 4.75
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("S13").Value = 1.9
$ws.Range("T13").Value = 1.95
$ws.Range("W13").Value = 3.2
$ws.Range("X13").Value = 1.36

# Row 15
$ws.Range("G15").Value = 2.2
$ws.Range("I15").Value = 3.2
$ws.Range("J15").Value = 2.88
$ws.Range("K15").Value = 2.05
$ws.Range("L15").Value = 4
$ws.Range("S15").Value = 2
$ws.Range("T15").Value = 1.8
$ws.Range("AC15").Value = 7
$ws.Range("AI15").Value = 9
$ws.Range("AP15").Value = 12

# Row 17
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 9
$ws.Range("S17").Value = 2.1
$ws.Range("T17").Value = 1.73

# Row 18
$ws.Range("M18").Value = 1.03
$ws.Range("N18").Value = 15
$ws.Range("W18").Value = 2.75
$ws.Range("X18").Value = 1.4

# Row 19
$ws.Range("N19").Value = 10

# Row 22
$ws.Range("N22").Value = 9
